$d = $word.ActiveDocument

# --- Body: bold run "QWREW" -> "QWR" ---
$body = $d.Content
$body.Find.Execute("QWREW", $true, $true, $false, $false, $false, $true, 1, $false, "QWR", 2)

# --- Header (primary header of section 1) ---
$h = $d.Sections(1).Headers(1)

# "DIRETORIA DE ENSINO REGIAO REW" -> "...QWER"  (standalone REW, not the REW inside QWREW)
$r1 = $h.Range
$r1.Find.Execute("REW", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 2)

# "QWREW - DEP." -> "QWR - DEP."
$r2 = $h.Range
$r2.Find.Execute("QWREW", $true, $true, $false, $false, $false, $true, 1, $false, "QWR", 2)

# "Rew, no Rew - Rew - Rew - Rew" -> "Qwer, no Qwer - Qwer - Qwer - Qwer"
$r3 = $h.Range
$r3.Find.Execute("Rew", $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 2)

# "CEP: rew ... Tel: rew" and "Email: rew" -> "qwer"
$r4 = $h.Range
$r4.Find.Execute("rew", $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 2)
